$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds scraped price strings that use "." as both a thousands
# separator and decimal point (e.g. "37.774.74", "14.40"). Several of the new
# values parse as plain numbers, which would make Excel silently coerce the
# cell to a Number and drop significant trailing zeros (e.g. "14.40" -> 14.4).
# Force those specific cells to Text before writing, then restore the default
# "Normal" style so no stray NumberFormat is left behind (matches source: all
# data cells use the default/no explicit style).

# Row 2
$ws.Range("D2").Value = "37.774.74"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").Value = "2.039.01"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.12%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0838"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.63%  "

# Row 11
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$ws.Range("D12").Value = "2.341.17"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.73%  "

# Row 14
$ws.Range("E14").Value = "  -0.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.95%  "

# Row 17
$ws.Range("D17").Value = "2.033.21"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18
$ws.Range("D18").Value = "37.717.02"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.42"
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.95%  "

# Row 21
$ws.Range("E21").Value = "  -0.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "

# Row 24
$ws.Range("E24").Value = "  +1.65%  "

# Row 25
$ws.Range("E25").Value = "  +2.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.45%  "

# Row 27
$ws.Range("E27").Value = "  +0.96%  "

# Row 28
$ws.Range("E28").Value = "  -0.49%  "

# Row 29
$ws.Range("E29").Value = "  -0.82%  "

# Row 30
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("E31").Value = "  -1.05%  "

# Row 32
$ws.Range("E32").Value = "  +8.64%  "

# Row 33
$ws.Range("E33").Value = "  -1.68%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "

# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.94%  "

# Row 37
$ws.Range("E37").Value = "  +3.44%  "

# Row 38
$ws.Range("E38").Value = "  +6.07%  "

# Row 39
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.84%  "

# Row 41
$ws.Range("D41").Value = "1.526.94"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.55%  "

# Row 43
$ws.Range("E43").Value = "  -1.28%  "

# Row 44
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.72%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.95%  "

# Row 47
$ws.Range("E47").Value = "  -0.54%  "

# Row 48
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("E49").Value = "  -0.37%  "

# Row 50
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("D51").Value = "2.230.44"
$ws.Range("E51").Value = "  +0.10%  "
